# R22 UAT 1 script update:
#  - D2 holds the shared-string value "POCOMM"; bulk upload expects "POCOMMS".
#  - Selection cursor should rest on D2 (was stuck on D11) after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the collection-type label in D2.
$ws.Range("D2").Value = "POCOMMS"

# Leave the active selection on the cell that was just edited.
$ws.Range("D2").Select()
